# Auto-generated: refresh market-price / profit columns (H..N) per scheduled runner diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 97.29412000000001
$ws.Range("I53").Value = 107.07143
$ws.Range("J53").Value = 51.666668
$ws.Range("K53").Value = 107.07143
$ws.Range("L53").Value = 51.666668
$ws.Range("M53").Value = 529.92857
$ws.Range("N53").Value = -1325.666668

$ws.Range("H98").Value = 760
$ws.Range("I98").Value = 766.4167
$ws.Range("K98").Value = 766.4167
$ws.Range("M98").Value = 731.5833

$ws.Range("H107").Value = 4480.421
$ws.Range("I107").Value = 1562.3334
$ws.Range("K107").Value = 1562.3334
$ws.Range("M107").Value = 357.6666

$ws.Range("H111").Value = 5903.696
$ws.Range("I111").Value = 5339.6313
$ws.Range("K111").Value = 16018.8939
$ws.Range("M111").Value = -12951.8939

$ws.Range("H122").Value = 760
$ws.Range("I122").Value = 766.4167
$ws.Range("K122").Value = 2299.2501
$ws.Range("M122").Value = 150.7498999999998

$ws.Range("H137").Value = 2238.2222
$ws.Range("I137").Value = 2205.5
$ws.Range("K137").Value = 6616.5
$ws.Range("M137").Value = -4066.5

$ws.Range("H138").Value = 3850.279
$ws.Range("J138").Value = 4041
$ws.Range("L138").Value = 12123
$ws.Range("N138").Value = -22403

$ws.Range("H139").Value = 133952.28
$ws.Range("J139").Value = 139533.2
$ws.Range("L139").Value = 139533.2
$ws.Range("N139").Value = -149813.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3471.5557
$ws.Range("I2").Value = 4290.6665
$ws.Range("J2").Value = 1833.3334
$ws.Range("K2").Value = 4290.6665
$ws.Range("L2").Value = 1833.3334
$ws.Range("M2").Value = -4177.6665
$ws.Range("N2").Value = -2059.3334

$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").Value = ""

$ws.Range("H38").Value = 20808.9
$ws.Range("I38").Value = 6454.3335
$ws.Range("J38").Value = 150000
$ws.Range("K38").Value = 6454.3335
$ws.Range("L38").Value = 150000
$ws.Range("M38").Value = -5987.3335
$ws.Range("N38").Value = -150934

$ws.Range("H97").Value = 1174
$ws.Range("I97").Value = 1258.28
$ws.Range("J97").Value = 120.5
$ws.Range("K97").Value = 1258.28
$ws.Range("L97").Value = 120.5
$ws.Range("M97").Value = -762.28
$ws.Range("N97").Value = -1112.5

$ws.Range("H116").Value = 3471.5557
$ws.Range("I116").Value = 4290.6665
$ws.Range("J116").Value = 1833.3334
$ws.Range("K116").Value = 4290.6665
$ws.Range("L116").Value = 1833.3334
$ws.Range("M116").Value = -1996.6665
$ws.Range("N116").Value = -6421.3334

$ws.Range("H122").Value = 2284.389
$ws.Range("I122").Value = 2002.3334
$ws.Range("K122").Value = 6007.0002
$ws.Range("M122").Value = -3557.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3471.5557
$ws.Range("I3").Value = 4290.6665
$ws.Range("J3").Value = 1833.3334
$ws.Range("K3").Value = 4290.6665
$ws.Range("L3").Value = 1833.3334
$ws.Range("M3").Value = -4176.6665
$ws.Range("N3").Value = -2061.3334

$ws.Range("H20").Value = 6718
$ws.Range("I20").Value = 2775
$ws.Range("J20").Value = 8032.3335
$ws.Range("K20").Value = 2775
$ws.Range("L20").Value = 8032.3335
$ws.Range("M20").Value = -2528
$ws.Range("N20").Value = -8526.333500000001

$ws.Range("H22").Value = 388310.88
$ws.Range("I22").Value = 523
$ws.Range("K22").Value = 523
$ws.Range("M22").Value = -350

$ws.Range("H94").Value = 1563.1428
$ws.Range("I94").Value = 1307.5238
$ws.Range("K94").Value = 1307.5238
$ws.Range("M94").Value = -856.5237999999999

$ws.Range("H99").Value = 5563.647
$ws.Range("I99").Value = 2334.9092
$ws.Range("J99").Value = 11483
$ws.Range("K99").Value = 2334.9092
$ws.Range("L99").Value = 11483
$ws.Range("M99").Value = -836.9092000000001
$ws.Range("N99").Value = -14479

$ws.Range("H105").Value = 3654.7
$ws.Range("J105").Value = 4666.3335
$ws.Range("L105").Value = 4666.3335
$ws.Range("N105").Value = -8160.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 29646.105
$ws.Range("J31").Value = 12959.8
$ws.Range("L31").Value = 12959.8
$ws.Range("N31").Value = -13549.8

$ws.Range("H34").Value = 29646.105
$ws.Range("J34").Value = 12959.8
$ws.Range("L34").Value = 12959.8
$ws.Range("N34").Value = -13363.8

$ws.Range("H62").Value = 7065.5557
$ws.Range("I62").Value = 7065.5557
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 7065.5557
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -6441.5557
$ws.Range("N62").Value = ""

$ws.Range("H65").Value = 7065.5557
$ws.Range("I65").Value = 7065.5557
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 35327.7785
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -32207.7785
$ws.Range("N65").Value = ""

$ws.Range("H107").Value = 591.9091
$ws.Range("I107").Value = 586.9643
$ws.Range("J107").Value = 619.6
$ws.Range("K107").Value = 586.9643
$ws.Range("L107").Value = 619.6
$ws.Range("M107").Value = 1333.0357
$ws.Range("N107").Value = -4459.6

$ws.Range("H122").Value = 1110.6538
$ws.Range("I122").Value = 1116.0952
$ws.Range("J122").Value = 1087.8
$ws.Range("K122").Value = 3348.2856
$ws.Range("L122").Value = 3263.4
$ws.Range("M122").Value = -898.2856000000002
$ws.Range("N122").Value = -8163.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 4105.4736
$ws.Range("J39").Value = 4111.3335
$ws.Range("L39").Value = 12334.0005
$ws.Range("N39").Value = -12922.0005

$ws.Range("H86").Value = 444.5
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").Value = ""

$ws.Range("H89").Value = 444.5
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").Value = ""

$ws.Range("H109").Value = 395.8
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = ""

$ws.Range("H125").Value = 3475
$ws.Range("I125").Value = 3475
$ws.Range("K125").Value = 10425
$ws.Range("M125").Value = -5505

$ws.Range("H140").Value = 14426.333
$ws.Range("I140").Value = 14426.333
$ws.Range("K140").Value = 43278.999
$ws.Range("M140").Value = -38098.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 232
$ws.Range("I2").Value = 197
$ws.Range("K2").Value = 197
$ws.Range("M2").Value = -84

$ws.Range("H102").Value = 2020.1471
$ws.Range("I102").Value = 2022.0968
$ws.Range("K102").Value = 2022.0968
$ws.Range("M102").Value = -400.0968

$ws.Range("H105").Value = 21500
$ws.Range("J105").Value = 21500
$ws.Range("L105").Value = 21500
$ws.Range("N105").Value = -28488

$ws.Range("H107").Value = 1340.8235
$ws.Range("I107").Value = 1182.1111
$ws.Range("J107").Value = 1519.375
$ws.Range("K107").Value = 1182.1111
$ws.Range("L107").Value = 1519.375
$ws.Range("M107").Value = 737.8888999999999
$ws.Range("N107").Value = -5359.375

$ws.Range("H124").Value = 28181.455
$ws.Range("J124").Value = 28181.455
$ws.Range("L124").Value = 28181.455
$ws.Range("N124").Value = -38001.455

$ws.Range("H126").Value = 24854
$ws.Range("I126").Value = 32255.545
$ws.Range("K126").Value = 96766.63499999999
$ws.Range("M126").Value = -94296.63499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4549.3
$ws.Range("I40").Value = 4651
$ws.Range("K40").Value = 4651
$ws.Range("M40").Value = -4515

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 31100
$ws.Range("J54").Value = 32200
$ws.Range("L54").Value = 32200
$ws.Range("N54").Value = -33240

$ws.Range("H126").Value = 4592.1904
$ws.Range("J126").Value = 4975
$ws.Range("L126").Value = 14925
$ws.Range("N126").Value = -19865

$ws.Range("H132").Value = 3088.3684
$ws.Range("I132").Value = 2704.389
$ws.Range("K132").Value = 8113.167
$ws.Range("M132").Value = -5583.167

$ws.Range("H133").Value = 89775
$ws.Range("J133").Value = 89775
$ws.Range("L133").Value = 89775
$ws.Range("N133").Value = -99895
